$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old 'NB' row (original row 8). This shifts the old SVM row
# (row 9) up to row 8, leaving 7 data rows (2..8) instead of 8 (2..9),
# matching the new dimension A1:L8.
$ws.Rows.Item(8).Delete()

# --- Header row (row 1) ---
# Columns C..G already carry the bordered/bold header style (s=1); extend
# that same formatting to the five new std columns (H..L) by copying the
# format from the existing "G1" header cell before filling in the text.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:L1").PasteSpecial(-4122) | Out-Null

$ws.Range("C1").Value = "One Year Alt mean"
$ws.Range("D1").Value = "One Year Alt std"
$ws.Range("E1").Value = "Two Year Alt mean"
$ws.Range("F1").Value = "Two Year Alt std"
$ws.Range("G1").Value = "Three Year Alt mean"
$ws.Range("H1").Value = "Three Year Alt std"
$ws.Range("I1").Value = "Five Year Alt mean"
$ws.Range("J1").Value = "Five Year Alt std"
$ws.Range("K1").Value = "Ten Year Alt mean"
$ws.Range("L1").Value = "Ten Year Alt std"

# --- Data rows (2..8) ---
# Row 2: LR
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.8250831845436508
$ws.Range("D2").Value = 0.0138980098013046
$ws.Range("E2").Value = 0.7999115016361236
$ws.Range("F2").Value = 0.008182131758107225
$ws.Range("G2").Value = 0.7788477071139324
$ws.Range("H2").Value = 0.01714833512599123
$ws.Range("I2").Value = 0.7597791704617372
$ws.Range("J2").Value = 0.03736497975085206
$ws.Range("K2").Value = 0.7380710692274287
$ws.Range("L2").Value = 0.03323787960698839

# Row 3: LDA
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.8304159480984952
$ws.Range("D3").Value = 0.01531690306900747
$ws.Range("E3").Value = 0.8059514657331628
$ws.Range("F3").Value = 0.01794649390305117
$ws.Range("G3").Value = 0.7855206638392469
$ws.Range("H3").Value = 0.01842389584805104
$ws.Range("I3").Value = 0.762627698032475
$ws.Range("J3").Value = 0.0361938755845554
$ws.Range("K3").Value = 0.7371322946399125
$ws.Range("L3").Value = 0.02434627599925425

# Row 4: KNN
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.813596638682734
$ws.Range("D4").Value = 0.01561174382783036
$ws.Range("E4").Value = 0.8126114365904955
$ws.Range("F4").Value = 0.02346061014179818
$ws.Range("G4").Value = 0.7968011194433338
$ws.Range("H4").Value = 0.01913415296520009
$ws.Range("I4").Value = 0.8003892066782949
$ws.Range("J4").Value = 0.02162452076189563
$ws.Range("K4").Value = 0.7810029239620618
$ws.Range("L4").Value = 0.02917181768783366

# Row 5: DTREE
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.7797780991180652
$ws.Range("D5").Value = 0.0246538586828926
$ws.Range("E5").Value = 0.7747980548851238
$ws.Range("F5").Value = 0.01691394957716253
$ws.Range("G5").Value = 0.7723093722228915
$ws.Range("H5").Value = 0.02405344085703648
$ws.Range("I5").Value = 0.7547810891083944
$ws.Range("J5").Value = 0.02903430617081509
$ws.Range("K5").Value = 0.7594763860428804
$ws.Range("L5").Value = 0.02965215656331961

# Row 6: RTREE
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8268175266283684
$ws.Range("D6").Value = 0.01161729592431021
$ws.Range("E6").Value = 0.805500666044374
$ws.Range("F6").Value = 0.01471296662425717
$ws.Range("G6").Value = 0.7850203372036102
$ws.Range("H6").Value = 0.01414790206283466
$ws.Range("I6").Value = 0.7624708089687131
$ws.Range("J6").Value = 0.03269955414404708
$ws.Range("K6").Value = 0.7400858145395569
$ws.Range("L6").Value = 0.03100621981639643

# Row 7: XTREE
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.8372464580411176
$ws.Range("D7").Value = 0.01843340251314643
$ws.Range("E7").Value = 0.8195601686493441
$ws.Range("F7").Value = 0.02189986229073317
$ws.Range("G7").Value = 0.8026616604770757
$ws.Range("H7").Value = 0.01999872103071963
$ws.Range("I7").Value = 0.8037101107565681
$ws.Range("J7").Value = 0.03050456643921642
$ws.Range("K7").Value = 0.7957760674428326
$ws.Range("L7").Value = 0.02108835906416959

# Row 8: SVM
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.8291553146432497
$ws.Range("D8").Value = 0.01013841243803716
$ws.Range("E8").Value = 0.8224202155615098
$ws.Range("F8").Value = 0.0146795624893133
$ws.Range("G8").Value = 0.8149535222722433
$ws.Range("H8").Value = 0.01423105559472435
$ws.Range("I8").Value = 0.8087109016781993
$ws.Range("J8").Value = 0.02754606637753372
$ws.Range("K8").Value = 0.7858349643772906
$ws.Range("L8").Value = 0.02527243457665334

